$d = $word.ActiveDocument

# Locate the unique word "line " that sits between "command " and
# "to schedule the importer" (it's in its own run) and remove it,
# along with the run break before/after it, by deleting just that
# span of text. This merges the three surrounding runs into one
# without disturbing the preceding "To schedule regular imports
# from GBIF" run.
$found = $d.Content.Duplicate
$found.Find.Execute("command line to schedule the importer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found.Find.Found) {
    # $found now spans "command line to schedule the importer";
    # narrow it down to just the "line " portion (7 chars in from
    # the start of "command ", i.e. right after "command ").
    $lineStart = $found.Start + 8
    $lineEnd = $lineStart + 5
    $toDelete = $d.Range($lineStart, $lineEnd)
    $toDelete.Delete()
}
